$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.598.69"
$ws.Range("E2").Value = "  -0.89%  "
Set-TextValue $ws.Range("D3") "3.446.71"
$ws.Range("E3").Value = "  -3.14%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("D5") "592.66"
$ws.Range("E5").Value = "  -1.80%  "
Set-TextValue $ws.Range("D6") "136.62"
$ws.Range("E6").Value = "  -6.98%  "
Set-TextValue $ws.Range("D7") "3.443.74"
$ws.Range("E7").Value = "  -3.37%  "
Set-TextValue $ws.Range("D8") "1.00"
$ws.Range("E8").Value = "  +0.10%  "
Set-TextValue $ws.Range("D9") "0.498"
$ws.Range("E9").Value = "  +2.21%  "
Set-TextValue $ws.Range("D10") "7.41"
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("E11").Value = "  -8.16%  "
Set-TextValue $ws.Range("D12") "0.379"
$ws.Range("E12").Value = "  -7.51%  "
Set-TextValue $ws.Range("D13") "4.025.33"
$ws.Range("E13").Value = "  -3.17%  "
Set-TextValue $ws.Range("D14") "0.0000182"
Set-TextValue $ws.Range("D15") "26.65"
$ws.Range("E15").Value = "  -9.00%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "65.552.26"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.115"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D18") "3.396.67"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("E19").Value = "  -10.55%  "
Set-TextValue $ws.Range("D20") "5.84"
$ws.Range("E20").Value = "  -6.64%  "
Set-TextValue $ws.Range("D21") "13.72"
$ws.Range("E21").Value = "  -6.92%  "
Set-TextValue $ws.Range("D22") "395.52"
$ws.Range("E22").Value = "  -5.84%  "
Set-TextValue $ws.Range("D23") "0.551"
$ws.Range("E23").Value = "  -8.94%  "
Set-TextValue $ws.Range("D24") "73.48"
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("E25").Value = "  -0.12%  "
Set-TextValue $ws.Range("D26") "3.585.85"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("E27").Value = "  -9.62%  "
Set-TextValue $ws.Range("D28") "0.996"
$ws.Range("E28").Value = "  -0.41%  "
Set-TextValue $ws.Range("D29") "7.23"
$ws.Range("E29").Value = "  -8.59%  "
$ws.Range("E30").Value = "  -8.81%  "
Set-TextValue $ws.Range("D31") "8.21"
$ws.Range("E31").Value = "  -11.06%  "
Set-TextValue $ws.Range("D32") "3.449.98"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -6.07%  "
Set-TextValue $ws.Range("D35") "23.07"
$ws.Range("E35").Value = "  -6.71%  "
Set-TextValue $ws.Range("D36") "171.42"
$ws.Range("E36").Value = "  -2.02%  "
Set-TextValue $ws.Range("D37") "6.98"
$ws.Range("E37").Value = "  -8.97%  "
$ws.Range("E38").Value = "  -10.86%  "
$ws.Range("E39").Value = "  -6.62%  "
$ws.Range("E40").Value = "  -9.93%  "
$ws.Range("E41").Value = "  -7.28%  "
Set-TextValue $ws.Range("D42") "0.826"
$ws.Range("E42").Value = "  -4.91%  "
Set-TextValue $ws.Range("D43") "43.60"
$ws.Range("E43").Value = "  -4.69%  "
Set-TextValue $ws.Range("D44") "0.999"
$ws.Range("E44").Value = "  +0.14%  "
Set-TextValue $ws.Range("D45") "4.44"
$ws.Range("E45").Value = "  -13.59%  "
$ws.Range("E46").Value = "  -11.33%  "
Set-TextValue $ws.Range("D47") "1.11"
$ws.Range("E47").Value = "  -1.41%  "
Set-TextValue $ws.Range("D48") "22.54"
$ws.Range("E48").Value = "  -3.38%  "
Set-TextValue $ws.Range("D49") "6.56"
$ws.Range("E49").Value = "  -7.46%  "
$ws.Range("E50").Value = "  -14.43%  "
Set-TextValue $ws.Range("D51") "2.205.12"
$ws.Range("E51").Value = "  -7.17%  "
